$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5699.3687
$ws.Range("I33").Value = 6412.375
$ws.Range("J33").Value = 1896.6666
$ws.Range("K33").Value = 6412.375
$ws.Range("L33").Value = 1896.6666
$ws.Range("M33").Value = -6183.375
$ws.Range("N33").Value = -2354.6666

$ws.Range("H41").Value = 366.58334
$ws.Range("I41").Value = 333.22223
$ws.Range("J41").Value = 466.66666
$ws.Range("K41").Value = 333.22223
$ws.Range("L41").Value = 466.66666
$ws.Range("M41").Value = 106.77777
$ws.Range("N41").Value = -1346.66666

$ws.Range("H87").Value = 19579.363
$ws.Range("J87").Value = 19579.363
$ws.Range("L87").Value = 19579.363
$ws.Range("N87").Value = -22075.363

$ws.Range("H90").Value = 19579.363
$ws.Range("J90").Value = 19579.363
$ws.Range("L90").Value = 58738.08900000001
$ws.Range("N90").Value = -71218.08900000001

$ws.Range("H132").Value = 38723.293
$ws.Range("I132").Value = 43712.633
$ws.Range("K132").Value = 131137.899
$ws.Range("M132").Value = -128607.899

$ws.Range("H138").Value = 15387859
$ws.Range("I138").Value = 1747.0385
$ws.Range("J138").Value = 25645268
$ws.Range("K138").Value = 5241.1155
$ws.Range("L138").Value = 76935804
$ws.Range("M138").Value = -101.1154999999999
$ws.Range("N138").Value = -76946084

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7120.0527
$ws.Range("I32").Value = 4705.0625
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 4705.0625
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -4418.0625
$ws.Range("N32").Value = -20574

$ws.Range("H45").Value = 2352.75
$ws.Range("I45").Value = 914.7778
$ws.Range("K45").Value = 914.7778
$ws.Range("M45").Value = -537.7778

$ws.Range("H74").Value = 2005.6923
$ws.Range("I74").Value = 2026
$ws.Range("J74").Value = 1938
$ws.Range("K74").Value = 2026
$ws.Range("L74").Value = 1938
$ws.Range("M74").Value = -1152
$ws.Range("N74").Value = -3686

$ws.Range("H77").Value = 2005.6923
$ws.Range("I77").Value = 2026
$ws.Range("J77").Value = 1938
$ws.Range("K77").Value = 10130
$ws.Range("L77").Value = 9690
$ws.Range("M77").Value = -5762
$ws.Range("N77").Value = -18426

$ws.Range("H122").Value = 1846.7046
$ws.Range("I122").Value = 1828.6061
$ws.Range("K122").Value = 5485.8183
$ws.Range("M122").Value = -3035.8183

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1462.375
$ws.Range("I134").Value = 928.42224
$ws.Range("J134").Value = 3646.7273
$ws.Range("K134").Value = 2785.26672
$ws.Range("L134").Value = 10940.1819
$ws.Range("M134").Value = -250.2667200000001
$ws.Range("N134").Value = -16010.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1210.4546
$ws.Range("I7").Value = 1873.5
$ws.Range("J7").Value = 414.8
$ws.Range("K7").Value = 1873.5
$ws.Range("L7").Value = 414.8
$ws.Range("M7").Value = -1760.5
$ws.Range("N7").Value = -640.8

$ws.Range("H58").Value = 2030.4706
$ws.Range("I58").Value = 1410.7273
$ws.Range("J58").Value = 3166.6667
$ws.Range("K58").Value = 1410.7273
$ws.Range("L58").Value = 3166.6667
$ws.Range("M58").Value = -1207.7273
$ws.Range("N58").Value = -3572.6667

$ws.Range("H105").Value = 202604
$ws.Range("I105").Value = 253005
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 253005
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -251258
$ws.Range("N105").Value = -4494

$ws.Range("H136").Value = 2030.4706
$ws.Range("I136").Value = 1410.7273
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 4232.1819
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -1682.1819
$ws.Range("N136").Value = -14600.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1595374.6
$ws.Range("I68").Value = 2961290.5
$ws.Range("J68").Value = 1806.1
$ws.Range("K68").Value = 8883871.5
$ws.Range("L68").Value = 5418.299999999999
$ws.Range("M68").Value = -8883060.5
$ws.Range("N68").Value = -7040.299999999999

$ws.Range("H71").Value = 1595374.6
$ws.Range("I71").Value = 2961290.5
$ws.Range("J71").Value = 1806.1
$ws.Range("K71").Value = 26651614.5
$ws.Range("L71").Value = 16254.9
$ws.Range("M71").Value = -26647558.5
$ws.Range("N71").Value = -24366.9

$ws.Range("H122").Value = 483.08572
$ws.Range("I122").Value = 313.66666
$ws.Range("J122").Value = 1499.6
$ws.Range("K122").Value = 2822.99994
$ws.Range("L122").Value = 13496.4
$ws.Range("M122").Value = -372.9999399999997
$ws.Range("N122").Value = -18396.4

$ws.Range("H137").Value = 1743.6857
$ws.Range("I137").Value = 1556.45
$ws.Range("J137").Value = 1993.3334
$ws.Range("K137").Value = 4669.35
$ws.Range("L137").Value = 5980.0002
$ws.Range("M137").Value = 430.6499999999996
$ws.Range("N137").Value = -16180.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14533
$ws.Range("J70").Value = 3932.3333
$ws.Range("L70").Value = 3932.3333
$ws.Range("N70").Value = -4472.3333

$ws.Range("H73").Value = 14533
$ws.Range("J73").Value = 3932.3333
$ws.Range("L73").Value = 3932.3333
$ws.Range("N73").Value = -5804.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58979.61
$ws.Range("I7").Value = 70048.53
$ws.Range("J7").Value = 3635
$ws.Range("K7").Value = 70048.53
$ws.Range("L7").Value = 3635
$ws.Range("M7").Value = -69936.53
$ws.Range("N7").Value = -3859

$ws.Range("H126").Value = 58979.61
$ws.Range("I126").Value = 70048.53
$ws.Range("J126").Value = 3635
$ws.Range("K126").Value = 210145.59
$ws.Range("L126").Value = 10905
$ws.Range("M126").Value = -207675.59
$ws.Range("N126").Value = -15845

$ws.Range("H132").Value = 3571.8071
$ws.Range("I132").Value = 3365.606
$ws.Range("J132").Value = 3855.3333
$ws.Range("K132").Value = 10096.818
$ws.Range("L132").Value = 11565.9999
$ws.Range("M132").Value = -7566.818000000001
$ws.Range("N132").Value = -16625.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3084.1292
$ws.Range("I136").Value = 474.44
$ws.Range("J136").Value = 4847.4326
$ws.Range("K136").Value = 1423.32
$ws.Range("L136").Value = 14542.2978
$ws.Range("M136").Value = 1126.68
$ws.Range("N136").Value = -19642.2978
